# Apollo "Digitale Kompetenzen" assessment workbook update
#  - adds courses / survey meta-data to row 2 (Duration, Publisher, Title)
#  - fixes J3 placeholder text to match J2 ("Bitte klicke an.")
#  - duplicates the "KldB" (AR) column value into the newly added
#    "DescriptionOfPartialQualification" (AT) column for every data row
#  - leaves the final selection on BI2 (matches the authored file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the stray single-space placeholder in J3 so it reuses the same
#     shared string as J2 ("Bitte klicke an.") -------------------------
$ws.Range("J3").Value = "Bitte klicke an."

# --- new survey / course meta-data on row 2 ----------------------------
$ws.Range("BF2").Value = 20
$ws.Range("BJ2").Value = "Selbsttest – Digitale Kompetenzen"
$ws.Range("BI2").Value = "Bertelsmann Stiftung"

# --- duplicate column AR ("KldB") into the new column AT
#     ("DescriptionOfPartialQualification") for data rows 2-56 ----------
for ($r = 2; $r -le 56; $r++) {
    $src = $ws.Range("AR$r").Value()
    if ($src -ne $null -and $src -ne "") {
        $ws.Range("AT$r").Value = $src
    }
}

# --- restore the selection left behind by the authoring session --------
$ws.Range("BI2").Select()
